# Apply scheduled-runner profit recalculations across Leve sheets (ALC/ARM/BSM/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
# Row 6: Days of Chunder / Antidote
$ws.Range("H6").Value = 73
$ws.Range("I6").Value = 51.444443
$ws.Range("K6").Value = 154.333329
$ws.Range("M6").Value = -42.33332899999999
# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 4999.6665
$ws.Range("I17").Value = 3000
$ws.Range("K17").Value = 9000
$ws.Range("M17").Value = -8832
# Row 29: Dripping with Venom / Weak Blinding Potion
$ws.Range("H29").Value = 2516.875
$ws.Range("I29").Value = 194.5
$ws.Range("J29").Value = 3291
$ws.Range("K29").Value = 583.5
$ws.Range("L29").Value = 9873
$ws.Range("M29").Value = -302.5
$ws.Range("N29").Value = -10435
# Row 115: 5-bell Energy / Competent Craftsman's Syrup
$ws.Range("H115").Value = 325
$ws.Range("I115").Value = 325
$ws.Range("K115").Value = 975
$ws.Range("M115").Value = 592
# Row 118: Crafty Concoctions / Commanding Craftsman's Syrup
$ws.Range("H118").Value = 1738.3334
$ws.Range("I118").Value = 1643.25
$ws.Range("K118").Value = 4929.75
$ws.Range("M118").Value = -3272.75
# Row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 3183.7144
$ws.Range("I129").Value = 1141.7142
$ws.Range("J129").Value = 4204.7144
$ws.Range("K129").Value = 3425.1426
$ws.Range("L129").Value = 12614.1432
$ws.Range("M129").Value = 1574.8574
$ws.Range("N129").Value = -22614.1432
# Row 131: Mindful Study / Grade 5 Tincture of Mind
$ws.Range("H131").Value = 1186.75
$ws.Range("J131").Value = 1200
$ws.Range("L131").Value = 3600
$ws.Range("N131").Value = -13680
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 4130.8184
$ws.Range("I132").Value = 4168.4443
$ws.Range("K132").Value = 12505.3329
$ws.Range("M132").Value = -9975.332900000001
# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 4484.3687
$ws.Range("I141").Value = 4372.1177
$ws.Range("J141").Value = 5438.5
$ws.Range("K141").Value = 13116.3531
$ws.Range("L141").Value = 16315.5
$ws.Range("M141").Value = -7936.3531
$ws.Range("N141").Value = -26675.5

$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 302.25
$ws.Range("I5").Value = 230.77777
$ws.Range("J5").Value = 516.6667
$ws.Range("K5").Value = 230.77777
$ws.Range("L5").Value = 516.6667
$ws.Range("M5").Value = -118.77777
$ws.Range("N5").Value = -740.6667
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 7045.1665
$ws.Range("I32").Value = 5817.8857
$ws.Range("K32").Value = 5817.8857
$ws.Range("M32").Value = -5530.8857
# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 2598
$ws.Range("I122").Value = 2598
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7794
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5344
$ws.Range("N122").ClearContents()
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 930
$ws.Range("I132").Value = 930
$ws.Range("K132").Value = 2790
$ws.Range("M132").Value = -260

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 302.25
$ws.Range("I4").Value = 230.77777
$ws.Range("J4").Value = 516.6667
$ws.Range("K4").Value = 230.77777
$ws.Range("L4").Value = 516.6667
$ws.Range("M4").Value = -115.77777
$ws.Range("N4").Value = -746.6667

$ws = $wb.Worksheets.Item("CUL")
# Row 39: Bloody Good Tart, This / Blood Currant Tart
$ws.Range("H39").Value = 1881.25
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -15588
# Row 55: Pagan Pastries / Pastry Fish
$ws.Range("H55").Value = 3366.3333
$ws.Range("I55").Value = 99
$ws.Range("K55").Value = 297
$ws.Range("M55").Value = -120
# Row 88: Don't Let It Fall Apart / Liver-cheese Sandwich
$ws.Range("H88").Value = 5000
$ws.Range("I88").Value = 5000
$ws.Range("K88").Value = 15000
$ws.Range("M88").Value = -14572
# Row 91: Better Come Back with a Sandwich (L) / Liver-cheese Sandwich
$ws.Range("H91").Value = 5000
$ws.Range("I91").Value = 5000
$ws.Range("K91").Value = 15000
$ws.Range("M91").Value = -13518
# Row 109: Cure for What Ails / Purple Carrot Juice
$ws.Range("H109").Value = 2000
$ws.Range("I109").Value = 2000
$ws.Range("K109").Value = 6000
$ws.Range("M109").Value = -4960
# Row 121: A Cookie for Your Troubles / Coffee Biscuit
$ws.Range("H121").Value = 976
$ws.Range("I121").Value = 895
$ws.Range("K121").Value = 2685
$ws.Range("M121").Value = -1375
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 2055.6667
$ws.Range("J131").Value = 2532.3076
$ws.Range("L131").Value = 7596.9228
$ws.Range("N131").Value = -17676.9228

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 3133.75
$ws.Range("I80").Value = 1642.5
$ws.Range("J80").Value = 4625
$ws.Range("K80").Value = 1642.5
$ws.Range("L80").Value = 4625
$ws.Range("M80").Value = -644.5
$ws.Range("N80").Value = -6621
# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 3133.75
$ws.Range("I83").Value = 1642.5
$ws.Range("J83").Value = 4625
$ws.Range("K83").Value = 8212.5
$ws.Range("L83").Value = 23125
$ws.Range("M83").Value = -3220.5
$ws.Range("N83").Value = -33109
# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 41733670
$ws.Range("I122").Value = 62550004
$ws.Range("J122").Value = 101008
$ws.Range("K122").Value = 187650012
$ws.Range("L122").Value = 303024
$ws.Range("M122").Value = -187647562
$ws.Range("N122").Value = -307924

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 7499
$ws.Range("I7").Value = 7498.5
$ws.Range("K7").Value = 7498.5
$ws.Range("M7").Value = -7386.5
# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 3359.8
$ws.Range("J46").Value = 3999.75
$ws.Range("L46").Value = 3999.75
$ws.Range("N46").Value = -4375.75
# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 714.5
$ws.Range("I61").Value = 714.5
$ws.Range("K61").Value = 714.5
$ws.Range("M61").Value = -512.5
# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 6666.6665
$ws.Range("I68").Value = 6666.6665
$ws.Range("K68").Value = 6666.6665
$ws.Range("M68").Value = -5917.6665
# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 6666.6665
$ws.Range("I71").Value = 6666.6665
$ws.Range("K71").Value = 33333.3325
$ws.Range("M71").Value = -29589.3325
# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 1616.9231
$ws.Range("I82").Value = 1553.3334
$ws.Range("K82").Value = 1553.3334
$ws.Range("M82").Value = -1192.3334
# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 1616.9231
$ws.Range("I85").Value = 1553.3334
$ws.Range("K85").Value = 1553.3334
$ws.Range("M85").Value = -305.3334
# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 714.5
$ws.Range("I113").Value = 714.5
$ws.Range("K113").Value = 714.5
$ws.Range("M113").Value = 1455.5
# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 7499
$ws.Range("I126").Value = 7498.5
$ws.Range("K126").Value = 22495.5
$ws.Range("M126").Value = -20025.5

$ws = $wb.Worksheets.Item("WVR")
# Row 100: Of Great Import / Kudzu Thread
$ws.Range("H100").Value = 711.6
$ws.Range("I100").Value = 639.875
$ws.Range("K100").Value = 1279.75
$ws.Range("M100").Value = -738.75
